$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1)
$ws.Range("A1").Value = "Semana"
$ws.Range("B1").Value = "Casos"
$ws.Range("C1").Value = "Temperatura"
$ws.Range("D1").Value = "Busquedas"

# Update column A (Semana) - sequential week numbers
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5

# Update column B (Casos) values
$ws.Range("B2").Value = 10
$ws.Range("B3").Value = 2
$ws.Range("B4").Value = 6
$ws.Range("B5").Value = 2
$ws.Range("B6").Value = 9

# Update selection to match the saved view state
$ws.Range("H6").Select()
